$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Text fix on "攻击链路分析页" (sheet 4): 登记严重度 -> 等级严重度
# ---------------------------------------------------------------------------
$wsAttack = $wb.Worksheets.Item(4)
$wsAttack.Range("B12").Value = "显示正在遭受攻击的端点IP，按顺序排列，且同时按照IP进行关联到各个其它群集；确认攻击事件红色连接（线的粗细按照等级严重度），疑似事件黄色连接"

# ---------------------------------------------------------------------------
# 2) Re-apply the existing centre alignment on "首页" (sheet 1) A1:C1 so the
#    cells keep pointing at the de-duplicated "center" style.
# ---------------------------------------------------------------------------
$wsHome = $wb.Worksheets.Item(1)
$wsHome.Range("A1:C1").HorizontalAlignment = -4108
$wsHome.Range("A1:C1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) Add the three new worksheets at the end of the workbook, in order.
# ---------------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsStorage.Name = "数据存储"

$wsBigData = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsBigData.Name = "大数据平台"

$wsML = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsML.Name = "机器学习平台"

# ---------------------------------------------------------------------------
# 4) Populate "数据存储"
# ---------------------------------------------------------------------------
$wsStorage.Range("A1").Value = "运维指标的数据处理方案："
$wsStorage.Range("A2").Value = "以JSON格式存储一份最初原始各指标内容，以后每次更改指标，只以JSON格式记录更新的指标，前端显示最新更新后的数据"
$wsStorage.Range("A3").Value = "考虑使用mongodb存储键值数据"
$wsStorage.Range("A5").Value = "对于实时日志数据处理方案："
$wsStorage.Range("A6").Value = "日志处理成JSON格式，存储到Hadoop中，再从Hadoop中抽取JSON数据，进行分析处理，如异常检测，攻击类型检测与确认"
$wsStorage.Range("A7").Value = "机器学习部分数据处理，可以考虑使用flink sql，和spark sql"
$wsStorage.Range("A9").Value = "对于攻击链路数据处理方案："
$wsStorage.Range("A10").Value = "可以用mongodb，和flink sql等"
$wsStorage.Range("A12").Value = "对于关联分析部分的数据处理方案："
$wsStorage.Range("A13").Value = "涉及图存储，考虑使用neo4j来存储关联关系"
$wsStorage.Columns.Item(1).ColumnWidth = 130.714285714286

# ---------------------------------------------------------------------------
# 5) Populate "大数据平台"
# ---------------------------------------------------------------------------
$wsBigData.Range("A1").Value = "spark"
$wsBigData.Range("A2").Value = "flink"
$wsBigData.Range("A3").Value = "hadoop"

# ---------------------------------------------------------------------------
# 6) Populate "机器学习平台"
# ---------------------------------------------------------------------------
$wsML.Range("A1").Value = "spark ml"
$wsML.Range("A2").Value = "flink ml"
$wsML.Range("A3").Value = "deeplearning4j"
$wsML.Range("A4").Value = "angel"
$wsML.Range("A5").Value = "spark graph"
$wsML.Columns.Item(1).ColumnWidth = 14.142857142857187

# ---------------------------------------------------------------------------
# 7) Restore/update the selections on every sheet to match the saved view.
# ---------------------------------------------------------------------------
$wsHome.Activate()
$wsHome.Range("B46").Select()

$wsAttack.Activate()
$wsAttack.Range("B16").Select()

$wsStorage.Activate()
$wsStorage.Range("A19").Select()

$wsML.Activate()
$wsML.Range("D11").Select()

$wsBigData.Activate()
$wsBigData.Range("K17").Select()

# The last activated sheet becomes the workbook's active tab, matching the
# target workbook.xml (activeTab points at "大数据平台").
